$d = $word.ActiveDocument

$map = @(
  @("681×9=6129", "471×8=3768"),
  @("171×3=513", "388×3=1164"),
  @("924×7=6468", "638×3=1914"),
  @("981×8=7848", "429×5=2145"),
  @("115×6=690", "324×7=2268"),
  @("301×9=2709", "474×8=3792"),
  @("518×7=3626", "671×6=4026"),
  @("333×5=1665", "262×8=2096"),
  @("245×8=1960", "324×8=2592"),
  @("464×4=1856", "117×5=585"),
  @("368×3=1104", "278×8=2224"),
  @("582×3=1746", "580×3=1740"),
  @("404×5=2020", "128×2=256"),
  @("313×4=1252", "788×5=3940"),
  @("577×4=2308", "784×8=6272"),
  @("709×5=3545", "526×4=2104"),
  @("374×7=2618", "804×3=2412"),
  @("450×5=2250", "275×5=1375"),
  @("585×2=1170", "486×2=972"),
  @("608×9=5472", "430×9=3870"),
  @("863×6=5178", "843×2=1686"),
  @("863×3=2589", "171×2=342"),
  @("393×6=2358", "843×3=2529"),
  @("858×3=2574", "418×7=2926"),
  @("992×3=2976", "626×6=3756")
)

foreach ($pair in $map) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
